$d = $word.ActiveDocument

$replacements = @(
    @("2023 年 7 月 25 日に Tailwind Traders と交渉されたサプライヤー契約条件は次のとおりです。", "2023 年 7 月 25 日に Tailwind Traders と交渉された供給者合意条件は次のとおりです。"),
    @("一律料金", "遅延料金"),
    @("1 か月あたり `$100", "1 か月あたり 1%"),
    @("早期支払い割引", "早期支払割引"),
    @("15 日以内に 5% 割引", "15 日以内 5% 割引"),
    @("契約は、いずれかの当事者が有効期限の少なくとも 90 日前に終了の書面による通知を行わない限り、さらに 3 年間自動的に更新されます。", "いずれかの当事者が有効期限の少なくとも 90 日前に書面による終了の通知を行わない限り、契約はさらに 3 年間自動的に更新されます。"),
    @("最小注文金額", "最小注文量"),
    @("20 時間/月", "1 か月あたり 20 ケース"),
    @("最大注文金額", "最大注文量"),
    @("1 か月あたり 100 リリース", "1 か月あたり 100 ケース")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
